# Update gh-pages output data (GZ con/event listing) to match newly
# generated scrape results at commit 456a3b4:
#  - refreshed "想去人数" (interest count, column F) for many events
#  - refreshed "最低票价" (lowest price, column G) for a couple of events
#  - the "街头霸王主题聚会" meetup was cancelled: title suffixed with
#    （取消） and its price marked as 不可售 (not available) on every
#    sheet that lists it ("展览" and "全部类型")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F3").Value = 1451
$ws.Range("F4").Value = 1110
$ws.Range("F5").Value = 526
$ws.Range("F8").Value = 678

$ws.Range("C10").Value = "广州·街头霸王主题聚会（取消）"
$ws.Range("G10").Value = "不可售"

$ws.Range("F11").Value = 89
$ws.Range("F12").Value = 219

$ws.Range("F14").Value = 2465
$ws.Range("G14").Value = 83.59999999999999

$ws.Range("F15").Value = 435
$ws.Range("F17").Value = 503
$ws.Range("F22").Value = 667
$ws.Range("F24").Value = 245
$ws.Range("F25").Value = 965
$ws.Range("F27").Value = 1577
$ws.Range("F28").Value = 315

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")

$ws.Range("F5").Value = 226

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F4").Value = 1451
$ws.Range("F5").Value = 1110
$ws.Range("F8").Value = 526
$ws.Range("F11").Value = 678

$ws.Range("C14").Value = "广州·街头霸王主题聚会（取消）"
$ws.Range("G14").Value = "不可售"

$ws.Range("F15").Value = 89
$ws.Range("F16").Value = 219

$ws.Range("F18").Value = 2465
$ws.Range("G18").Value = 83.59999999999999

$ws.Range("F19").Value = 226
$ws.Range("F20").Value = 435
$ws.Range("F22").Value = 503
$ws.Range("F31").Value = 667
$ws.Range("F37").Value = 245
$ws.Range("F38").Value = 965
$ws.Range("F40").Value = 1577
$ws.Range("F41").Value = 315
